# river update May 2024
# Append three new observation rows (42-44) for "Manganui o te Ao at Ashworth"
# sampled 2023-04-05 (Excel serial date 45021).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$siteName = "Manganui o te Ao at Ashworth"
$sampleDate = 45021

$rows = @(
    @{ Row = 42; Param = "ASPM (Macroinvertebrate Average Score Per Metric)"; Value = "0.391" },
    @{ Row = 43; Param = "MCI (Macroinvertebrate Community Index)";            Value = "104.76" },
    @{ Row = 44; Param = "QMCI (Quantitative Macroinvertebrate Community Index)"; Value = "4.079" }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $siteName
    $ws.Cells.Item($row, 2).Value = $r.Param

    # date time column keeps the same date-time style/number format as the row above
    $ws.Cells.Item($row, 3).Value = $sampleDate
    $ws.Cells.Item($row, 3).Style = $ws.Cells.Item($row - 1, 3).Style
    $ws.Cells.Item($row, 3).NumberFormat = $ws.Cells.Item($row - 1, 3).NumberFormat

    # Value column stores numeric-looking text as TEXT (not a number) -
    # force text entry via a "@" format, then drop back to the Normal
    # style so no extra formatting is left on the cell.
    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $r.Value
    $ws.Cells.Item($row, 4).Style = "Normal"

    # Project / Method / Unit / pH columns are present but blank
    $ws.Cells.Item($row, 5).Formula = '=""'
    $ws.Cells.Item($row, 6).Formula = '=""'
    $ws.Cells.Item($row, 7).Formula = '=""'

    $ws.Cells.Item($row, 8).Value = 200

    $ws.Cells.Item($row, 9).Formula = '=""'
}
